$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 keeps the same visible text ("Number of Wells "); no edit required -
# its shared-string index will renumber automatically once unused strings
# are pruned below.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Row 6: "Average Daily Withdrawal (MGD)" - numeric values, keeps style s="1"
# (B6:F6 already carry style s="1" in the source file, so only the cell
# values themselves need to change)
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Average Daily Withdrawal (MGD)"
$ws.Range("B6").Value = 0.03
$ws.Range("C6").Value = 0.03
$ws.Range("D6").Value = 0.03
$ws.Range("E6").Value = 0.03
$ws.Range("F6").Value = 0.04

# ---------------------------------------------------------------------------
# Row 7: "Average Daily Withdrawal (gpd) " - numeric values, new style
# (font ArialMT 10pt, number format #,##0), label loses its old bold style
# ---------------------------------------------------------------------------
$ws.Range("A7").Style = "Normal"
$ws.Range("A7").Value = "Average Daily Withdrawal (gpd) "

# B7 already carries style s="1" (ArialMT 10pt) in the source file; adding a
# #,##0 number format on top of it produces the new combined style.
$ws.Range("B7").Value = 29000
$ws.Range("B7").NumberFormat = "#,##0"

$ws.Range("B7").Copy()
$ws.Range("C7:F7").PasteSpecial(-4122)
$ws.Range("C7").Value = 27800
$ws.Range("D7").Value = 25500
$ws.Range("E7").Value = 25700
$ws.Range("F7").Value = 40000

# ---------------------------------------------------------------------------
# Row 8: "Design Capacity- Max Daily " - numeric values, style s="1" kept,
# row is now taller (25pt)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Design Capacity- Max Daily "
$ws.Range("B8").Value = 0.07
$ws.Range("C8").Value = 0.04
$ws.Range("D8").Value = 0.04
$ws.Range("E8").Value = 0.04
$ws.Range("F8").Value = 0.12
$ws.Rows.Item(8).RowHeight = 25

# ---------------------------------------------------------------------------
# Row 9: "Design Capacity- Max Daily " (duplicate label) - numeric values
# with the new #,##0 style, row is also taller (25pt)
# ---------------------------------------------------------------------------
$ws.Range("B11").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "Design Capacity- Max Daily "

$ws.Range("B7").Copy()
$ws.Range("B9:F9").PasteSpecial(-4122)
$ws.Range("B9").Value = 68400
$ws.Range("C9").Value = 43200
$ws.Range("D9").Value = 36000
$ws.Range("E9").Value = 36000
$ws.Range("F9").Value = 115200
$ws.Rows.Item(9).RowHeight = 25

# ---------------------------------------------------------------------------
# Row 10: "System Permitted Capacity: " - text values (same text previously
# shown in rows 7/8), style s="1" kept, row is also taller (25pt)
# ---------------------------------------------------------------------------
$ws.Range("B11").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)
$ws.Range("A10").Value = "System Permitted Capacity: "
$ws.Range("B10").Value = "0.07 MGD (68,400 gpd) "
$ws.Range("C10").Value = "0.04 MGD (43,200 gpd) "
$ws.Range("D10").Value = "0.04 MGD (36,000 gpd) "
$ws.Range("E10").Value = "0.04 MGD (36,000 gpd) "
$ws.Range("F10").Value = "0.12 MGD (115,200 gpd) "
$ws.Rows.Item(10).RowHeight = 25

# ---------------------------------------------------------------------------
# sheet view: zoom + selected cell
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 139
$ws.Range("C10").Select()
